$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row (B11): right-answer mark value 3 -> 5
$ws.Range("B11").Value = 5

# Update "Total" row (B12): total correct score 72 -> 120
$ws.Range("B12").Value = 120

# Update "Total" row (E12): corr/total marks text "71/84" -> "120/140"
$ws.Range("E12").Value = "120/140"
